$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(51, 8).Value = 23895.8
$ws.Cells.Item(51, 9).Value = 35833.332
$ws.Cells.Item(51, 10).Value = 5989.5
$ws.Cells.Item(51, 11).Value = 35833.332
$ws.Cells.Item(51, 12).Value = 5989.5
$ws.Cells.Item(51, 13).Value = -35349.332
$ws.Cells.Item(51, 14).Value = -6957.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(53, 8).Value = 278.32144
$ws.Cells.Item(53, 9).Value = 227.07692
$ws.Cells.Item(53, 10).Value = 322.73334
$ws.Cells.Item(53, 11).Value = 227.07692
$ws.Cells.Item(53, 12).Value = 322.73334
$ws.Cells.Item(53, 13).Value = 409.92308
$ws.Cells.Item(53, 14).Value = -1596.73334

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(129, 8).Value = 292330.97
$ws.Cells.Item(129, 10).Value = 302845.34
$ws.Cells.Item(129, 12).Value = 908536.02
$ws.Cells.Item(129, 14).Value = -918536.02

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 2526.487
$ws.Cells.Item(138, 9).Value = 2215.6
$ws.Cells.Item(138, 10).Value = 2600.508
$ws.Cells.Item(138, 11).Value = 6646.799999999999
$ws.Cells.Item(138, 12).Value = 7801.523999999999
$ws.Cells.Item(138, 13).Value = -1506.799999999999
$ws.Cells.Item(138, 14).Value = -18081.524

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 31198.594
$ws.Cells.Item(32, 9).Value = 5594.06
$ws.Cells.Item(32, 11).Value = 5594.06
$ws.Cells.Item(32, 13).Value = -5307.06

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1581.5714
$ws.Cells.Item(61, 9).Value = 1498.2693
$ws.Cells.Item(61, 11).Value = 1498.2693
$ws.Cells.Item(61, 13).Value = -1286.2693

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(103, 8).Value = 0
$ws.Cells.Item(103, 10).Value = 0
$ws.Cells.Item(103, 12).Value = 0
$ws.Cells.Item(103, 14).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 77077800
$ws.Cells.Item(110, 9).Value = 83500900
$ws.Cells.Item(110, 10).Value = 600
$ws.Cells.Item(110, 11).Value = 83500900
$ws.Cells.Item(110, 12).Value = 600
$ws.Cells.Item(110, 13).Value = -83498855
$ws.Cells.Item(110, 14).Value = -4690

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 2143.9048
$ws.Cells.Item(122, 9).Value = 2113.8
$ws.Cells.Item(122, 10).Value = 2219.1667
$ws.Cells.Item(122, 11).Value = 6341.400000000001
$ws.Cells.Item(122, 12).Value = 6657.500100000001
$ws.Cells.Item(122, 13).Value = -3891.400000000001
$ws.Cells.Item(122, 14).Value = -11557.5001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 1581.5714
$ws.Cells.Item(136, 9).Value = 1498.2693
$ws.Cells.Item(136, 11).Value = 4494.8079
$ws.Cells.Item(136, 13).Value = -1944.8079

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(45, 8).Value = 14166.667
$ws.Cells.Item(45, 9).Value = 12500
$ws.Cells.Item(45, 11).Value = 12500
$ws.Cells.Item(45, 13).Value = -11907

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 7973.575
$ws.Cells.Item(58, 9).Value = 1383.3043
$ws.Cells.Item(58, 11).Value = 1383.3043
$ws.Cells.Item(58, 13).Value = -1180.3043

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 1241.3549
$ws.Cells.Item(134, 9).Value = 1153.1538
$ws.Cells.Item(134, 11).Value = 3459.4614
$ws.Cells.Item(134, 13).Value = -924.4614000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 7973.575
$ws.Cells.Item(136, 9).Value = 1383.3043
$ws.Cells.Item(136, 11).Value = 4149.9129
$ws.Cells.Item(136, 13).Value = -1599.9129

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 8521.77
$ws.Cells.Item(5, 9).Value = 745.1177
$ws.Cells.Item(5, 10).Value = 23211
$ws.Cells.Item(5, 11).Value = 2235.3531
$ws.Cells.Item(5, 12).Value = 69633
$ws.Cells.Item(5, 13).Value = -2123.3531
$ws.Cells.Item(5, 14).Value = -69857

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 204.6
$ws.Cells.Item(34, 9).Value = 130
$ws.Cells.Item(34, 10).Value = 503
$ws.Cells.Item(34, 11).Value = 390
$ws.Cells.Item(34, 12).Value = 1509
$ws.Cells.Item(34, 13).Value = -306
$ws.Cells.Item(34, 14).Value = -1677

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(55, 8).Value = 10497.952
$ws.Cells.Item(55, 9).Value = 100000
$ws.Cells.Item(55, 10).Value = 6022.85
$ws.Cells.Item(55, 11).Value = 300000
$ws.Cells.Item(55, 12).Value = 18068.55
$ws.Cells.Item(55, 13).Value = -299823
$ws.Cells.Item(55, 14).Value = -18422.55

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(121, 8).Value = 7245.0454
$ws.Cells.Item(121, 9).Value = 6889.857
$ws.Cells.Item(121, 10).Value = 7410.8
$ws.Cells.Item(121, 11).Value = 20669.571
$ws.Cells.Item(121, 12).Value = 22232.4
$ws.Cells.Item(121, 13).Value = -19359.571
$ws.Cells.Item(121, 14).Value = -24852.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 5232.2383
$ws.Cells.Item(122, 9).Value = 384.2857
$ws.Cells.Item(122, 10).Value = 14928.143
$ws.Cells.Item(122, 11).Value = 3458.5713
$ws.Cells.Item(122, 12).Value = 134353.287
$ws.Cells.Item(122, 13).Value = -1008.5713
$ws.Cells.Item(122, 14).Value = -139253.287

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 769.09
$ws.Cells.Item(131, 10).Value = 863.8
$ws.Cells.Item(131, 12).Value = 2591.4
$ws.Cells.Item(131, 14).Value = -12671.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(135, 8).Value = 8521.77
$ws.Cells.Item(135, 9).Value = 745.1177
$ws.Cells.Item(135, 10).Value = 23211
$ws.Cells.Item(135, 11).Value = 6706.0593
$ws.Cells.Item(135, 12).Value = 208899
$ws.Cells.Item(135, 13).Value = -4171.0593
$ws.Cells.Item(135, 14).Value = -213969

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(9, 8).Value = 1303.5
$ws.Cells.Item(9, 9).Value = 1404.6666
$ws.Cells.Item(9, 10).Value = 1000
$ws.Cells.Item(9, 11).Value = 1404.6666
$ws.Cells.Item(9, 12).Value = 1000
$ws.Cells.Item(9, 13).Value = -1234.6666
$ws.Cells.Item(9, 14).Value = -1340

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 131264.62
$ws.Cells.Item(70, 9).Value = 226956.11
$ws.Cells.Item(70, 10).Value = 8232.714
$ws.Cells.Item(70, 11).Value = 226956.11
$ws.Cells.Item(70, 12).Value = 8232.714
$ws.Cells.Item(70, 13).Value = -226686.11
$ws.Cells.Item(70, 14).Value = -8772.714

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 131264.62
$ws.Cells.Item(73, 9).Value = 226956.11
$ws.Cells.Item(73, 10).Value = 8232.714
$ws.Cells.Item(73, 11).Value = 226956.11
$ws.Cells.Item(73, 12).Value = 8232.714
$ws.Cells.Item(73, 13).Value = -226020.11
$ws.Cells.Item(73, 14).Value = -10104.714

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 628.5714
$ws.Cells.Item(122, 9).Value = 712.5
$ws.Cells.Item(122, 10).Value = 516.6667
$ws.Cells.Item(122, 11).Value = 2137.5
$ws.Cells.Item(122, 12).Value = 1550.0001
$ws.Cells.Item(122, 13).Value = 312.5
$ws.Cells.Item(122, 14).Value = -6450.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 57566.832
$ws.Cells.Item(40, 9).Value = 112744.89
$ws.Cells.Item(40, 10).Value = 2388.7778
$ws.Cells.Item(40, 11).Value = 112744.89
$ws.Cells.Item(40, 12).Value = 2388.7778
$ws.Cells.Item(40, 13).Value = -112608.89
$ws.Cells.Item(40, 14).Value = -2660.7778

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 844346.8
$ws.Cells.Item(46, 9).Value = 500
$ws.Cells.Item(46, 10).Value = 1447094.6
$ws.Cells.Item(46, 11).Value = 500
$ws.Cells.Item(46, 12).Value = 1447094.6
$ws.Cells.Item(46, 13).Value = -312
$ws.Cells.Item(46, 14).Value = -1447470.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(64, 8).Value = 20000
$ws.Cells.Item(64, 10).Value = 20000
$ws.Cells.Item(64, 12).Value = 20000
$ws.Cells.Item(64, 14).Value = -20450

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(67, 8).Value = 20000
$ws.Cells.Item(67, 10).Value = 20000
$ws.Cells.Item(67, 12).Value = 20000
$ws.Cells.Item(67, 14).Value = -21560

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(81, 8).Value = 37165.5
$ws.Cells.Item(81, 10).Value = 37165.5
$ws.Cells.Item(81, 12).Value = 37165.5
$ws.Cells.Item(81, 14).Value = -39161.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(84, 8).Value = 37165.5
$ws.Cells.Item(84, 10).Value = 37165.5
$ws.Cells.Item(84, 12).Value = 111496.5
$ws.Cells.Item(84, 14).Value = -121480.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(105, 8).Value = 49990
$ws.Cells.Item(105, 10).Value = 49990
$ws.Cells.Item(105, 12).Value = 49990
$ws.Cells.Item(105, 14).Value = -56978

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(76, 8).Value = 173
$ws.Cells.Item(76, 10).Value = 173
$ws.Cells.Item(76, 12).Value = 173
$ws.Cells.Item(76, 14).Value = -803

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(79, 8).Value = 173
$ws.Cells.Item(79, 10).Value = 173
$ws.Cells.Item(79, 12).Value = 173
$ws.Cells.Item(79, 14).Value = -2357

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(97, 8).Value = 32379
$ws.Cells.Item(97, 10).Value = 32379
$ws.Cells.Item(97, 12).Value = 32379
$ws.Cells.Item(97, 14).Value = -34361

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2216
$ws.Cells.Item(122, 9).Value = 1201.3334
$ws.Cells.Item(122, 10).Value = 2554.2222
$ws.Cells.Item(122, 11).Value = 3604.0002
$ws.Cells.Item(122, 12).Value = 7662.6666
$ws.Cells.Item(122, 13).Value = -1154.0002
$ws.Cells.Item(122, 14).Value = -12562.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 1527.4286
$ws.Cells.Item(126, 9).Value = 1776
$ws.Cells.Item(126, 10).Value = 1080
$ws.Cells.Item(126, 11).Value = 5328
$ws.Cells.Item(126, 12).Value = 3240
$ws.Cells.Item(126, 13).Value = -2858
$ws.Cells.Item(126, 14).Value = -8180
